# Auto-push with fireball accuracy report [2025-04-16 04:55 PM]
# Update win/loss counts and recompute totals/accuracy on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "ATS Accuracy" ---
$ws1 = $wb.Worksheets.Item(1)

# Row 2: 5 loss, 86 -> 87 win
$ws1.Range("C2").Value = 87
$ws1.Range("D2").Value = 92
$ws1.Range("E2").Value = 94.59999999999999

# Row 3: 4 -> 6 loss, 66 win
$ws1.Range("B3").Value = 6
$ws1.Range("D3").Value = 72
$ws1.Range("E3").Value = 91.7

# Row 5: 5 -> 6 loss, 8 win
$ws1.Range("B5").Value = 6
$ws1.Range("D5").Value = 14
$ws1.Range("E5").Value = 57.1

# Row 6: 2 -> 3 loss, 5 win
$ws1.Range("B6").Value = 3
$ws1.Range("D6").Value = 8
$ws1.Range("E6").Value = 62.5

# --- Sheet 2: "Total Accuracy" ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2: 7 -> 8 loss, 71 win
$ws2.Range("B2").Value = 8
$ws2.Range("D2").Value = 79
$ws2.Range("E2").Value = 89.90000000000001

# Row 3: 2 -> 3 loss, 71 win
$ws2.Range("B3").Value = 3
$ws2.Range("D3").Value = 74
$ws2.Range("E3").Value = 95.90000000000001

# Row 4: 5 -> 6 loss, 16 -> 17 win
$ws2.Range("B4").Value = 6
$ws2.Range("C4").Value = 17
$ws2.Range("D4").Value = 23
$ws2.Range("E4").Value = 73.90000000000001

# Row 6: 1 loss, 3 -> 4 win
$ws2.Range("C6").Value = 4
$ws2.Range("D6").Value = 5
$ws2.Range("E6").Value = 80
